# "Fixed ICDC breed all testcases"
#
# Updates the cartQuery (column D) Cypher text on the "startup" sheet for
# the three data tabs (Cases/Samples/Files), rewriting the stats query so it
# walks program->study->case->demographic via OPTIONAL MATCHes and returns
# Programs/Studies/Cases/Samples/Case Files/Study Files counts instead of
# the old number_of_files/number_of_sample/number_of_cases/number_of_study
# shape. Also moves the saved cursor/selection from C2 to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCartQuery = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
    "OPTIONAL MATCH (samp:sample)-->(c)`n" +
    "OPTIONAL MATCH (diag:diagnosis)-->(c)`n" +
    "OPTIONAL MATCH (f:file)-[*]->(c)`n" +
    "OPTIONAL MATCH (sf:file)-->(s)`n" +
    "WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`n" +
    "WHERE demo.breed IN ['Yorkshire Terrier']`n" +
    "RETURN  `n" +
    "    count(distinct p) AS Programs,`n" +
    "    count(distinct s) AS Studies,`n" +
    "    count(distinct c) AS Cases,`n" +
    "    count(distinct samp) AS Samples,`n" +
    "    count(distinct f) AS ``Case Files``,`n" +
    "    count(distinct sf) AS ``Study Files``"

$ws.Range("D2").Value2 = $newCartQuery
$ws.Range("D3").Value2 = $newCartQuery
$ws.Range("D4").Value2 = $newCartQuery

# Move the saved selection/cursor down to C4 (from C2).
[void]$ws.Range("C4").Select()
